$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7154053167776624
$ws.Range("C2").Value = 0.1770292306925967
$ws.Range("E2").Value = 0.8077042195962889
$ws.Range("F2").Value = 1.98919547651596
$ws.Range("G2").Value = 0.2316775927200041
$ws.Range("H2").Value = 0.3994847428385029
$ws.Range("J2").Value = 0.02288687042903703
$ws.Range("N2").Value = 0.8629917256815958
$ws.Range("O2").Value = 1.176644007352508

$ws.Range("B3").Value = 0.6271257959401737
$ws.Range("C3").Value = 0.1585651896770912
$ws.Range("E3").Value = 0.7828031414453989
$ws.Range("F3").Value = 1.953259537971689
$ws.Range("G3").Value = 0.2300315429398125
$ws.Range("H3").Value = 0.402453762122299
$ws.Range("J3").Value = 0.02294979719200541
$ws.Range("N3").Value = 0.8575136044288314
$ws.Range("O3").Value = 1.179113091472573

$ws.Range("B4").Value = 0.5727779517782494
$ws.Range("C4").Value = 0.1471688742461481
$ws.Range("E4").Value = 0.7679053570022631
$ws.Range("F4").Value = 1.932451238492092
$ws.Range("G4").Value = 0.2292751221624698
$ws.Range("H4").Value = 0.4045180687315906
$ws.Range("J4").Value = 0.02300880137955552
$ws.Range("N4").Value = 0.8544993169160904
$ws.Range("O4").Value = 1.181664448710094

$ws.Range("B5").Value = 0.5505960729733204
$ws.Range("C5").Value = 0.1425102232924758
$ws.Range("E5").Value = 0.761932780931744
$ws.Range("F5").Value = 1.924287061634971
$ws.Range("G5").Value = 0.2290305354598701
$ws.Range("H5").Value = 0.4054199217971899
$ws.Range("J5").Value = 0.02303795105700601
$ws.Range("N5").Value = 0.8533591653662711
$ws.Range("O5").Value = 1.182963796545678

$ws.Range("B6").Value = 0.5469107374038913
$ws.Range("C6").Value = 0.1417357889787354
$ws.Range("E6").Value = 0.7609469798131414
$ws.Range("F6").Value = 1.922950437108568
$ws.Range("G6").Value = 0.2289937589268689
$ws.Range("H6").Value = 0.4055733343482899
$ws.Range("J6").Value = 0.02304309918715397
$ws.Range("N6").Value = 0.8531751815471296
$ws.Range("O6").Value = 1.183195214130308

$ws.Range("B7").Value = 0.5724789376971842
$ws.Range("C7").Value = 0.1471061045284046
$ws.Range("E7").Value = 0.767824410463831
$ws.Range("F7").Value = 1.932339857390161
$ws.Range("G7").Value = 0.2292715661850906
$ws.Range("H7").Value = 0.4045299860100755
$ws.Range("J7").Value = 0.02300917385375989
$ws.Range("N7").Value = 0.8544835828540016
$ws.Range("O7").Value = 1.181680921782473

$ws.Range("B8").Value = 0.6849972589871527
$ws.Range("C8").Value = 0.1706753492635755
$ws.Range("E8").Value = 0.7990370013544634
$ws.Range("F8").Value = 1.976543605612846
$ws.Range("G8").Value = 0.2310570964286711
$ws.Range("H8").Value = 0.4004583451435764
$ws.Range("J8").Value = 0.02290433048027651
$ws.Range("N8").Value = 0.8610305828508871
$ws.Range("O8").Value = 1.177280044782364

$ws.Range("B9").Value = 0.9044491442074332
$ws.Range("C9").Value = 0.2164103471588135
$ws.Range("E9").Value = 0.8633592775955066
$ws.Range("F9").Value = 2.073230667863442
$ws.Range("G9").Value = 0.2365889381735542
$ws.Range("H9").Value = 0.3943910497924179
$ws.Range("J9").Value = 0.02286110591126445
$ws.Range("N9").Value = 0.8766266995901901
$ws.Range("O9").Value = 1.176897715459233

$ws.Range("B10").Value = 1.064893091331385
$ws.Range("C10").Value = 0.2497015699131566
$ws.Range("E10").Value = 0.9125319068949409
$ws.Range("F10").Value = 2.150422525640693
$ws.Range("G10").Value = 0.241909776185139
$ws.Range("H10").Value = 0.3911058357879398
$ws.Range("J10").Value = 0.02292944514429962
$ws.Range("N10").Value = 0.8897498550976337
$ws.Range("O10").Value = 1.181692783589995

$ws.Range("B11").Value = 1.137700887438825
$ws.Range("C11").Value = 0.2647762359464423
$ws.Range("E11").Value = 0.9353215826878625
$ws.Range("F11").Value = 2.186888856881041
$ws.Range("G11").Value = 0.2446073002337243
$ws.Range("H11").Value = 0.3898667043965105
$ws.Range("J11").Value = 0.02298251090132197
$ws.Range("N11").Value = 0.8960781867138365
$ws.Range("O11").Value = 1.184986821095066

$ws.Range("B12").Value = 1.16524427444898
$ws.Range("C12").Value = 0.2704742797189965
$ws.Range("E12").Value = 0.9440121463445479
$ws.Range("F12").Value = 2.200892877802517
$ws.Range("G12").Value = 0.2456689326386652
$ws.Range("H12").Value = 0.3894342576428187
$ws.Range("J12").Value = 0.02300578473407811
$ws.Range("N12").Value = 0.8985258166780454
$ws.Range("O12").Value = 1.186394997375146

$ws.Range("B13").Value = 1.159313555023402
$ws.Range("C13").Value = 0.2692475729319312
$ws.Range("E13").Value = 0.9421377808210423
$ws.Range("F13").Value = 2.197868177090385
$ws.Range("G13").Value = 0.245438501172643
$ws.Range("H13").Value = 0.389525755640264
$ws.Range("J13").Value = 0.02300063061584723
$ws.Range("N13").Value = 0.8979964025749325
$ws.Range("O13").Value = 1.186084557490574

$ws.Range("B14").Value = 1.139967455713247
$ws.Range("C14").Value = 0.2652452281732565
$ws.Range("E14").Value = 0.9360353453452035
$ws.Range("F14").Value = 2.188037063810782
$ws.Range("G14").Value = 0.244693835221895
$ws.Range("H14").Value = 0.3898303891093065
$ws.Range("J14").Value = 0.0229843618352632
$ws.Range("N14").Value = 0.8962785298841993
$ws.Range("O14").Value = 1.185099445524685

$ws.Range("B15").Value = 1.128113802397536
$ws.Range("C15").Value = 0.2627923099774989
$ws.Range("E15").Value = 0.9323053233154326
$ws.Range("F15").Value = 2.182040639954948
$ws.Range("G15").Value = 0.2442429419251368
$ws.Range("H15").Value = 0.3900217785145088
$ws.Range("J15").Value = 0.02297481128681866
$ws.Range("N15").Value = 0.8952329449862049
$ws.Range("O15").Value = 1.184516999642199

$ws.Range("B16").Value = 1.060131191661128
$ws.Range("C16").Value = 0.2487149699502709
$ws.Range("E16").Value = 0.9110510238204483
$ws.Range("F16").Value = 2.148066616789748
$ws.Range("G16").Value = 0.2417390871764269
$ws.Range("H16").Value = 0.3911919597403255
$ws.Range("J16").Value = 0.02292642101274467
$ws.Range("N16").Value = 0.8893434702001741
$ws.Range("O16").Value = 1.181499968792508

$ws.Range("B17").Value = 1.018379078043893
$ws.Range("C17").Value = 0.2400608568329972
$ws.Range("E17").Value = 0.8981200448727407
$ws.Range("F17").Value = 2.127571259614371
$ws.Range("G17").Value = 0.2402742253854626
$ws.Range("H17").Value = 0.3919752732322621
$ws.Range("J17").Value = 0.02290237743917345
$ws.Range("N17").Value = 0.885822062884273
$ws.Range("O17").Value = 1.179934676174014

$ws.Range("B18").Value = 0.9943476179778941
$ws.Range("C18").Value = 0.2350767043695328
$ws.Range("E18").Value = 0.8907220894651005
$ws.Range("F18").Value = 2.115910027996847
$ws.Range("G18").Value = 0.2394577348678837
$ws.Range("H18").Value = 0.3924498430391878
$ws.Range("J18").Value = 0.02289061576399831
$ws.Range("N18").Value = 0.8838304095637852
$ws.Range("O18").Value = 1.179139050459781

$ws.Range("B19").Value = 0.9862081511008682
$ws.Range("C19").Value = 0.2333880452768256
$ws.Range("E19").Value = 0.8882240675863358
$ws.Range("F19").Value = 2.111983553590335
$ws.Range("G19").Value = 0.2391857517704636
$ws.Range("H19").Value = 0.3926146490457114
$ws.Range("J19").Value = 0.02288698800111888
$ws.Range("N19").Value = 0.8831618796414347
$ws.Range("O19").Value = 1.178887621777847

$ws.Range("B20").Value = 1.022825406917832
$ws.Range("C20").Value = 0.2409827811434013
$ws.Range("E20").Value = 0.8994924703360141
$ws.Range("F20").Value = 2.129739860259406
$ws.Range("G20").Value = 0.2404274632181966
$ws.Range("H20").Value = 0.3918894008000677
$ws.Range("J20").Value = 0.02290472279982936
$ws.Range("N20").Value = 0.8861934303507155
$ws.Range("O20").Value = 1.1800904636072

$ws.Range("B21").Value = 1.145650628992655
$ws.Range("C21").Value = 0.2664210993029883
$ws.Range("E21").Value = 0.9378261328919706
$ws.Range("F21").Value = 2.190919401998769
$ws.Range("G21").Value = 0.2449114699311679
$ws.Range("H21").Value = 0.3897399119752549
$ws.Range("J21").Value = 0.0229890539460591
$ws.Range("N21").Value = 0.8967817230177531
$ws.Range("O21").Value = 1.18538442658496

$ws.Range("B22").Value = 1.225763878981695
$ws.Range("C22").Value = 0.2829856840838545
$ws.Range("E22").Value = 0.963232714947253
$ws.Range("F22").Value = 2.232040774884211
$ws.Range("G22").Value = 0.2480761242624396
$ws.Range("H22").Value = 0.3885495426321484
$ws.Range("J22").Value = 0.02306270739433103
$ws.Range("N22").Value = 0.9040002322341962
$ws.Range("O22").Value = 1.189781941610534

$ws.Range("B23").Value = 1.183021112850838
$ws.Range("C23").Value = 0.2741505460632538
$ws.Range("E23").Value = 0.9496403823119834
$ws.Range("F23").Value = 2.209989259784066
$ws.Range("G23").Value = 0.2463655701283614
$ws.Range("H23").Value = 0.3891652203658964
$ws.Range("J23").Value = 0.02302169477701455
$ws.Range("N23").Value = 0.9001203774456741
$ws.Range("O23").Value = 1.187348855107473

$ws.Range("B24").Value = 1.020815307524686
$ws.Range("C24").Value = 0.2405660065065547
$ws.Range("E24").Value = 0.8988718839173089
$ws.Range("F24").Value = 2.128759056484796
$ws.Range("G24").Value = 0.2403581044312517
$ws.Range("H24").Value = 0.3919281482475441
$ws.Range("J24").Value = 0.02290365604326183
$ws.Range("N24").Value = 0.8860254327648107
$ws.Range("O24").Value = 1.180019707332491

$ws.Range("B25").Value = 0.8452156840631346
$ws.Range("C25").Value = 0.2040912224809404
$ws.Range("E25").Value = 0.8456232885102395
$ws.Range("F25").Value = 2.045997381604039
$ws.Range("G25").Value = 0.2348731122920285
$ws.Range("H25").Value = 0.3958267888412763
$ws.Range("J25").Value = 0.02285530151278792
$ws.Range("N25").Value = 0.8721140438009769
$ws.Range("O25").Value = 1.176113209035009

